# Auto-generated edit script: refreshes cached Universalis market-price
# columns (currentAveragePrice*, LevePriceNQ/HQ, LeveProfitNQ/HQ) across
# the Seraph_Profits workbook sheets, per the scheduled-runner commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 95.875
$ws.Range("I9").Value = 102.42857
$ws.Range("K9").Value = 102.42857
$ws.Range("M9").Value = 66.57143000000001
$ws.Range("H17").Value = 1330
$ws.Range("J17").Value = 1343.35
$ws.Range("L17").Value = 4030.05
$ws.Range("N17").Value = -4366.049999999999
$ws.Range("H55").Value = 1193.5625
$ws.Range("J55").Value = 1699.75
$ws.Range("L55").Value = 1699.75
$ws.Range("N55").Value = -2127.75
$ws.Range("H62").Value = 5166
$ws.Range("I62").Value = 4749
$ws.Range("K62").Value = 4749
$ws.Range("M62").Value = -4125
$ws.Range("H65").Value = 5166
$ws.Range("I65").Value = 4749
$ws.Range("K65").Value = 23745
$ws.Range("M65").Value = -20625
$ws.Range("H80").Value = 1006.2727
$ws.Range("I80").Value = 1073.2307
$ws.Range("J80").Value = 909.55554
$ws.Range("K80").Value = 3219.6921
$ws.Range("L80").Value = 2728.66662
$ws.Range("M80").Value = -2221.6921
$ws.Range("N80").Value = -4724.66662
$ws.Range("H83").Value = 1006.2727
$ws.Range("I83").Value = 1073.2307
$ws.Range("J83").Value = 909.55554
$ws.Range("K83").Value = 9659.076300000001
$ws.Range("L83").Value = 8185.99986
$ws.Range("M83").Value = -4667.076300000001
$ws.Range("N83").Value = -18169.99986
$ws.Range("H106").Value = 34614.69
$ws.Range("I106").Value = 35271.91
$ws.Range("K106").Value = 35271.91
$ws.Range("M106").Value = -34640.91
$ws.Range("H132").Value = 935.4828
$ws.Range("I132").Value = 893.7037
$ws.Range("K132").Value = 2681.1111
$ws.Range("M132").Value = -151.1111000000001
$ws.Range("H137").Value = 1632.4667
$ws.Range("I137").Value = 1387.6666
$ws.Range("J137").Value = 1999.6666
$ws.Range("K137").Value = 4162.9998
$ws.Range("L137").Value = 5998.9998
$ws.Range("M137").Value = -1612.9998
$ws.Range("N137").Value = -11098.9998
$ws.Range("H138").Value = 4190.079
$ws.Range("J138").Value = 4978.933
$ws.Range("L138").Value = 14936.799
$ws.Range("N138").Value = -25216.799
$ws.Range("H141").Value = 2079.4473
$ws.Range("I141").Value = 1381.4572
$ws.Range("K141").Value = 4144.3716
$ws.Range("M141").Value = 1035.6284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 3416.6667
$ws.Range("J15").Value = 4125
$ws.Range("L15").Value = 4125
$ws.Range("N15").Value = -4825
$ws.Range("H102").Value = 3999.5
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 12
$ws.Range("I132").Value = 12
$ws.Range("K132").Value = 36
$ws.Range("M132").Value = 2494
$ws.Range("H105").Value = 2825.2666
$ws.Range("I105").Value = 2567.6924
$ws.Range("K105").Value = 2567.6924
$ws.Range("M105").Value = -820.6923999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 126.411766
$ws.Range("I7").Value = 50.454544
$ws.Range("J7").Value = 265.66666
$ws.Range("K7").Value = 50.454544
$ws.Range("L7").Value = 265.66666
$ws.Range("M7").Value = 62.545456
$ws.Range("N7").Value = -491.66666
$ws.Range("H50").Value = 32970.832
$ws.Range("J50").Value = 32970.832
$ws.Range("L50").Value = 32970.832
$ws.Range("N50").Value = -34220.832
$ws.Range("H51").Value = 43387.2
$ws.Range("J51").Value = 43387.2
$ws.Range("L51").Value = 43387.2
$ws.Range("N51").Value = -44859.2
$ws.Range("H59").Value = 59999
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 59999
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 59999
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -62289
$ws.Range("H60").Value = 10863.625
$ws.Range("J60").Value = 49979
$ws.Range("L60").Value = 49979
$ws.Range("N60").Value = -51001
$ws.Range("H61").Value = 43387.2
$ws.Range("J61").Value = 43387.2
$ws.Range("L61").Value = 43387.2
$ws.Range("N61").Value = -44083.2
$ws.Range("H105").Value = 1045.5
$ws.Range("I105").Value = 729.7646999999999
$ws.Range("K105").Value = 729.7646999999999
$ws.Range("M105").Value = 1017.2353
$ws.Range("H107").Value = 896.8333
$ws.Range("J107").Value = 829.5
$ws.Range("L107").Value = 829.5
$ws.Range("N107").Value = -4669.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 34857520
$ws.Range("I4").Value = 53043708
$ws.Range("J4").Value = 660.75
$ws.Range("K4").Value = 159131124
$ws.Range("L4").Value = 1982.25
$ws.Range("M4").Value = -159131012
$ws.Range("N4").Value = -2206.25
$ws.Range("H6").Value = 517.6667
$ws.Range("I6").Value = 750
$ws.Range("J6").Value = 53
$ws.Range("K6").Value = 2250
$ws.Range("L6").Value = 159
$ws.Range("M6").Value = -2137
$ws.Range("N6").Value = -385
$ws.Range("H134").Value = 1842.6666
$ws.Range("I134").Value = 1842.6666
$ws.Range("K134").Value = 5527.9998
$ws.Range("M134").Value = -457.9997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 10333.333
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 10333.333
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 10333.333
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -10823.333
$ws.Range("H41").Value = 627.2
$ws.Range("I41").Value = 378.66666
$ws.Range("K41").Value = 378.66666
$ws.Range("M41").Value = -23.66665999999998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 136.6
$ws.Range("I9").Value = 136.6
$ws.Range("K9").Value = 136.6
$ws.Range("M9").Value = 87.40000000000001
$ws.Range("H11").Value = 1399.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1399.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 1399.5
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -1679.5
$ws.Range("H13").Value = 5556500
$ws.Range("I13").Value = 5883294
$ws.Range("J13").Value = 999
$ws.Range("K13").Value = 5883294
$ws.Range("L13").Value = 999
$ws.Range("M13").Value = -5883154
$ws.Range("N13").Value = -1279
$ws.Range("H14").Value = 19320.5
$ws.Range("J14").Value = 19320.5
$ws.Range("L14").Value = 19320.5
$ws.Range("N14").Value = -19664.5
$ws.Range("H30").Value = 1953.2
$ws.Range("I30").Value = 1897.25
$ws.Range("K30").Value = 1897.25
$ws.Range("M30").Value = -1789.25
$ws.Range("H42").Value = 7016007.5
$ws.Range("J42").Value = 7016007.5
$ws.Range("L42").Value = 7016007.5
$ws.Range("N42").Value = -7017133.5
$ws.Range("H46").Value = 4148.9
$ws.Range("I46").Value = 3000
$ws.Range("K46").Value = 3000
$ws.Range("M46").Value = -2812
$ws.Range("H49").Value = 7016007.5
$ws.Range("J49").Value = 7016007.5
$ws.Range("L49").Value = 7016007.5
$ws.Range("N49").Value = -7016301.5
$ws.Range("H55").Value = 778.46155
$ws.Range("I55").Value = 626.8333
$ws.Range("J55").Value = 908.4286
$ws.Range("K55").Value = 626.8333
$ws.Range("L55").Value = 908.4286
$ws.Range("M55").Value = -453.8333
$ws.Range("N55").Value = -1254.4286
$ws.Range("H93").Value = 1807.3334
$ws.Range("I93").Value = 1549.1666
$ws.Range("J93").Value = 2065.5
$ws.Range("K93").Value = 1549.1666
$ws.Range("L93").Value = 2065.5
$ws.Range("M93").Value = -301.1666
$ws.Range("N93").Value = -4561.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 5000
$ws.Range("I14").Value = 5000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -4832
$ws.Range("N14").ClearContents()
$ws.Range("H62").Value = 7235.0586
$ws.Range("I62").Value = 4750
$ws.Range("K62").Value = 4750
$ws.Range("M62").Value = -4126
$ws.Range("H65").Value = 7235.0586
$ws.Range("I65").Value = 4750
$ws.Range("K65").Value = 23750
$ws.Range("M65").Value = -20630
$ws.Range("H92").Value = 16767
$ws.Range("J92").Value = 16767
$ws.Range("L92").Value = 16767

Write-Host "Applied 210 edits across 8 sheets"